$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking price cells as Text so Excel does not coerce them to numbers
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D11", "D15", "D17", "D18", "D20", "D22", "D23", "D24", "D25", "D30", "D31", "D34", "D36", "D37", "D40", "D41", "D43", "D44", "D49", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated price values (Price column, D)
$ws.Range("D2").Value = "64.065.53"
$ws.Range("D3").Value = "3.151.70"
$ws.Range("D5").Value = "603.80"
$ws.Range("D6").Value = "143.94"
$ws.Range("D8").Value = "3.144.72"
$ws.Range("D11").Value = "5.39"
$ws.Range("D15").Value = "3.671.03"
$ws.Range("D17").Value = "64.109.92"
$ws.Range("D18").Value = "3.149.87"
$ws.Range("D20").Value = "491.06"
$ws.Range("D22").Value = "0.712"
$ws.Range("D23").Value = "7.67"
$ws.Range("D24").Value = "88.00"
$ws.Range("D25").Value = "13.31"
$ws.Range("D30").Value = "2.06"
$ws.Range("D31").Value = "27.74"
$ws.Range("D34").Value = "2.66"
$ws.Range("D36").Value = "6.04"
$ws.Range("D37").Value = "52.70"
$ws.Range("D40").Value = "434.99"
$ws.Range("D41").Value = "0.0396"
$ws.Range("D43").Value = "8.32"
$ws.Range("D44").Value = "2.934.81"
$ws.Range("D49").Value = "25.91"
$ws.Range("D51").Value = "120.43"

# Restore the General number format now that the text values are safely stored
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "General"
}

# Apply remaining text updates (Coin, Link, Volume(1h) columns)
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("E6").Value = "  -3.05%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("E23").Value = "  -3.90%  "
$ws.Range("E24").Value = "  +3.98%  "
$ws.Range("E25").Value = "  -3.57%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("E28").Value = "  -4.19%  "
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("E31").Value = "  +4.38%  "
$ws.Range("E32").Value = "  -5.13%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("E38").Value = "  -4.59%  "
$ws.Range("E39").Value = "  -6.64%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("E40").Value = "  -5.79%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("E45").Value = "  -3.65%  "
$ws.Range("E46").Value = "  -5.77%  "
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  +0.04%  "
